$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.678.54"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "1.638.81"
$ws.Range("E3").Value = "  -0.60%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'212.46"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("D6").Value = "'0.524"
$ws.Range("E6").Value = "  -1.60%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "'23.12"
$ws.Range("E8").Value = "  -2.11%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("E10").Value = "  -0.10%  "

$ws.Range("D11").Value = "'0.0894"
$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("D12").Value = "1.870.91"
$ws.Range("E12").Value = "  -0.52%  "

$ws.Range("D13").Value = "1.641.06"
$ws.Range("E13").Value = "  -0.50%  "

$ws.Range("E14").Value = "  +0.25%  "

$ws.Range("E15").Value = "  -5.11%  "

$ws.Range("D16").Value = "'64.69"
$ws.Range("E16").Value = "  +0.22%  "

$ws.Range("D17").Value = "27.662.86"
$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("D18").Value = "'230.84"
$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("D19").Value = "'7.71"
$ws.Range("E19").Value = "  +1.88%  "

$ws.Range("E20").Value = "  -0.49%  "

$ws.Range("E22").Value = "  -0.71%  "

$ws.Range("D23").Value = "'10.23"
$ws.Range("E23").Value = "  +4.52%  "

$ws.Range("D24").Value = "'2.03"
$ws.Range("E24").Value = "  +0.69%  "

$ws.Range("D25").Value = "'151.37"

$ws.Range("E26").Value = "  -1.01%  "

$ws.Range("E27").Value = "  -1.54%  "

$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("E30").Value = "  +0.14%  "

$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  -0.57%  "

$ws.Range("D33").Value = "1.458.06"
$ws.Range("E33").Value = "  +2.31%  "

$ws.Range("E34").Value = "  -1.15%  "

$ws.Range("E35").Value = "  -1.98%  "

$ws.Range("E36").Value = "  -0.12%  "

$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("E38").Value = "  -1.31%  "

$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("E40").Value = "  +9.41%  "

$ws.Range("E41").Value = "  +7.15%  "

$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'1.02"
$ws.Range("E42").Value = "  -0.81%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").Value = "'5.61"
$ws.Range("E44").Value = "  +1.37%  "

$ws.Range("E45").Value = "  -0.56%  "

$ws.Range("E46").Value = "  -0.42%  "

$ws.Range("D47").Value = "1.781.06"
$ws.Range("E47").Value = "  -0.51%  "

$ws.Range("D48").Value = "'1.75"
$ws.Range("E48").Value = "  +3.61%  "

$ws.Range("D49").Value = "'86.81"
$ws.Range("E49").Value = "  -1.60%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0993"
$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.79"
$ws.Range("E51").Value = "  +0.36%  "

# Clear quote-prefix styling so forced-text cells keep the default style
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
